$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add missing weekly Trazab_indic value for the existing last row (F72) ---
# Copy formatting from the existing row above so the new/edited cells keep
# the same number formats / styles as the rest of the table.
$ws.Cells.Item(72, 6).Value = 40.6

# --- Append the new data row (row 73) for 2020-12-24 ---
# Copy the date/territory formatting (style) from row 72 down into row 73
# first, then fill in the values (mirrors what Excel does when you type a
# new row right after the last one and it auto-extends formatting).
$ws.Range("A72:B72").Copy()
$ws.Range("A73:B73").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(73, 1).Value = 44189
$ws.Cells.Item(73, 2).Value = "Andalucía"
$ws.Cells.Item(73, 3).Value = 91202
$ws.Cells.Item(73, 4).Value = 8

# Update selection / view to point at the newly entered cell, matching
# where a user's cursor would land after typing the last value.
[void]$ws.Range("F73").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
